$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column S (2022 data) values per row, and whether the row needs the
# bold "category header" styling (mirrors the existing bold rows 7,10,13,
# 16,19,22,25,28,31 that already use a bold font for column A).
$values = @{
    3  = "2022"
    4  = "1.9210869108320343"
    5  = "1.020872301352429"
    6  = "2.8415499553180767"
    7  = "1.5924017665043597"
    8  = "2.5011433798307796"
    9  = "0.70098698968147144"
    10 = "2.2312343573160249"
    11 = "2.4764236727529938"
    12 = "1.9888745417939038"
    13 = "1.3057776932131271"
    14 = "2.6056788910230639"
    15 = "0"
    16 = "0.65058422463372112"
    17 = "0.65686622262510019"
    18 = "0.64442124527961442"
    19 = "2.5553368555544047"
    20 = "1.807815324711445"
    21 = "3.2928586128833093"
    22 = "1.8387963974300983"
    23 = "2.2260807622100529"
    24 = "1.4582467499325562"
    25 = "1.2245886088767601"
    26 = "1.3105423773238725"
    27 = "1.1375464261135158"
    28 = "2.4791112740241377"
    29 = "2.4279584268771761"
    30 = "2.5408788313520994"
    31 = "1.1238322680339958"
    32 = "0.57553956834532372"
    33 = "1.6467682173734046"
}

$boldRows = @(7, 10, 13, 16, 19, 22, 25, 28, 31)

for ($row = 3; $row -le 33; $row++) {
    # Column R (18) already carries the correct number formatting / borders
    # for this row; clone that formatting into the new column S (19) cell
    # before writing the value so the new column visually matches the rest
    # of the 2007-2021 series.
    $ws.Range($ws.Cells.Item($row, 18), $ws.Cells.Item($row, 18)).Copy()
    $target = $ws.Cells.Item($row, 19)
    $target.PasteSpecial(-4122)

    if ($boldRows -contains $row) {
        $target.Font.Bold = $true
    }

    $target.Value = [double]$values[$row]
}

$excel.CutCopyMode = $false

# Move the active selection the way the source workbook's author left it
# after adding the new 2022 column.
$ws.Range("T3").Select()
